# Apply the diff: rename "Velázquez" -> "Tienda Velázquez" (sheet name and column A),
# and fill the previously-empty column B (Nombre_TPV) with "BAR" for rows 2-5
# and "SERVIDOR TIENDA" for rows 6-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab name)
$ws.Name = "Tienda Velázquez"

# Update column A (Tienda) values for data rows 2-11
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "Tienda Velázquez"
}

# Update column B (Nombre_TPV) values
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 2).Value = "BAR"
}
for ($r = 6; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = "SERVIDOR TIENDA"
}
